# Doing Updates for Financials
# Insert a new "period" column (D) in front of the existing data (old D:K
# shifts right to E:L) and populate it with the latest fiscal-year figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column at D; existing D:K data moves to E:L.
$ws.Columns("D").Insert()

# 2. Copy the number formatting from the (now shifted) E column onto the new
#    D column so the new cells pick up the same date/number styles as the
#    rest of the row, without minting new style entries. Done in separate
#    blocks that line up exactly with the sheet's existing row ranges so we
#    don't manufacture cells in rows that never had D:K data (37, 79) or in
#    row numbers that don't exist at all (36, 78).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3. Fill in the new column D with the newest period's values.

# -- Income Statement --
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 550200
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 77200
$ws.Range("D18").Value = 472900
$ws.Range("D20").Value = -202700
$ws.Range("D21").Value = 284600
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 270200
$ws.Range("D24").Value = 50900
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 219300
$ws.Range("D27").Value = 219300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 300
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 202700
$ws.Range("D33").Value = 219600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 219600

# -- Balance Sheet --
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 324100
$ws.Range("D42").Value = 323600
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 151800
$ws.Range("D49").Value = 55800
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 17144000
$ws.Range("D57").Value = 8300
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 10600
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 10600
$ws.Range("D62").Value = 49600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 15875800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 1641300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1268200
$ws.Range("D77").Value = 0

# -- Cash Flow Statement --
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 219600
$ws.Range("D83").Value = 14400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 317900
$ws.Range("D91").Value = -35300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -74700
$ws.Range("D96").Value = -98500
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -165100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 78100
